# cryptos.xlsx refresh -- GitHub Actions scheduled update
# Updates the Price (D) and Volume(1h) (E) columns for each coin row with a
# fresh snapshot, and reflects that Hedera/ImmutableX swapped rank (rows
# 33/34), so their Coin/Link/Price/Volume are exchanged along with new
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33/34 rank swap: ImmutableX now outranks Hedera.
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"

# Per-row Price (D) / Volume(1h) (E) updates. Price is formatted text (the
# feed uses "." as both thousands separator and decimal point, e.g.
# "24.659.07", so it can never be a real Excel number) -- entries that would
# otherwise be auto-parsed as a plain number get a leading apostrophe to
# force text, then the style is reset back to Normal so no explicit
# NumberFormat/quote-prefix styling is left behind on the cell.
$updates = @(
    @{ Row = 2;  D = "24.659.07";  E = "  -0.17%  " }
    @{ Row = 3;  D = "1.685.99";   E = "  -0.99%  " }
    @{ Row = 4;  D = "1.003";      E = "  +0.80%  " }
    @{ Row = 5;  D = "315.32";     E = "  +0.52%  " }
    @{ Row = 6;  E = "  +0.62%  " }
    @{ Row = 7;  E = "  -0.81%  " }
    @{ Row = 8;  E = "  -0.81%  " }
    @{ Row = 9;  D = "1.002";      E = "  +0.67%  " }
    @{ Row = 10; E = "  -2.21%  " }
    @{ Row = 11; D = "53.28";      E = "  +0.05%  " }
    @{ Row = 12; D = "0.08792";    E = "  +0.02%  " }
    @{ Row = 13; D = "7.226";      E = "  -1.38%  " }
    @{ Row = 14; D = "23.40";      E = "  +0.41%  " }
    @{ Row = 15; D = "8.034";      E = "  +7.08%  " }
    @{ Row = 16; E = "  -1.15%  " }
    @{ Row = 17; D = "1.692.64";   E = "  -0.51%  " }
    @{ Row = 18; D = "99.47";      E = "  -1.56%  " }
    @{ Row = 19; D = "0.07004";    E = "  -1.36%  " }
    @{ Row = 20; D = "19.42";      E = "  -0.31%  " }
    @{ Row = 21; D = "6.980";      E = "  +3.44%  " }
    @{ Row = 22; D = "1.006";      E = "  +1.03%  " }
    @{ Row = 23; E = "  +0.51%  " }
    @{ Row = 24; D = "24.672.11";  E = "  -0.04%  " }
    @{ Row = 25; D = "3.277";      E = "  +10.10%  " }
    @{ Row = 26; D = "2.361";      E = "  +2.34%  " }
    @{ Row = 27; D = "22.64";      E = "  +0.86%  " }
    @{ Row = 28; D = "162.38";     E = "  +2.28%  " }
    @{ Row = 29; D = "5.185";      E = "  +0.89%  " }
    @{ Row = 30; D = "134.77";     E = "  +1.04%  " }
    @{ Row = 31; D = "7.564";      E = "  +1.43%  " }
    @{ Row = 32; D = "1.878.42";   E = "  -0.47%  " }
    @{ Row = 33; D = "1.057";      E = "  -3.04%  " }
    @{ Row = 34; D = "0.08531";    E = "  -1.67%  " }
    @{ Row = 35; D = "7.143";      E = "  -3.17%  " }
    @{ Row = 36; D = "11.13";      E = "  -0.34%  " }
    @{ Row = 37; D = "0.2720";     E = "  -0.13%  " }
    @{ Row = 38; D = "1.883";      E = "  -3.29%  " }
    @{ Row = 39; D = "14.34";      E = "  -3.14%  " }
    @{ Row = 40; D = "0.09150";    E = "  +1.84%  " }
    @{ Row = 41; D = "0.02709";    E = "  -2.10%  " }
    @{ Row = 42; D = "1.458";      E = "  -1.75%  " }
    @{ Row = 43; D = "0.7586";     E = "  -0.80%  " }
    @{ Row = 44; D = "15.89";      E = "  +2.71%  " }
    @{ Row = 45; D = "2.587";      E = "  +5.16%  " }
    @{ Row = 46; D = "0.7109";     E = "  -1.57%  " }
    @{ Row = 47; D = "4.218";      E = "  +1.20%  " }
    @{ Row = 48; E = "  +0.75%  " }
    @{ Row = 49; D = "139.21";     E = "  -1.61%  " }
    @{ Row = 50; D = "1.308";      E = "  -1.07%  " }
    @{ Row = 51; D = "0.07958";    E = "  -0.58%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $priceCell = $ws.Range("D" + $u.Row)
        $priceText = $u.D
        $isNumericLooking = $priceText -match '^[+-]?[0-9]+(\.[0-9]+)?$'

        if ($isNumericLooking) {
            # Force text storage (quote-prefix), then strip the quote-prefix
            # style back off so the cell keeps its original default styling.
            $priceCell.Value = "'" + $priceText
            $priceCell.Style = "Normal"
        } else {
            $priceCell.Value = $priceText
        }
    }

    $ws.Range("E" + $u.Row).Value = $u.E
}
